$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.173.37"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.803.44"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4668"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +24.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3714"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07693"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.386"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.409"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.25%  "
$ws.Range("D16").Value = "1.800.03"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001096"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06733"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.436"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.29%  "
$ws.Range("D23").Value = "28.153.28"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.411"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.401"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "2.006.34"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.271"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.042"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09652"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.929"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2259"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02378"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06413"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6743"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.273"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.521"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.235"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.132"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("E44").Value = "  +3.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9990"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6202"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.841"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.071"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.190"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07154"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.62%  "
